# Apply updated fee/reward figures per commit: "refactor to address issues resulting in bloated fees"
$wb = $excel.ActiveWorkbook

# --- Positions sheet: feeGrowthInsidePeriod0/1_formatted (columns F/G) ---
$wsPositions = $wb.Worksheets.Item("Positions")
$wsPositions.Range("G10").Value = "'596.81"
$wsPositions.Range("F18").Value = "'0.003209645145741389"
$wsPositions.Range("G18").Value = "'5681.33"
$wsPositions.Range("G20").Value = "'1131.44"
$wsPositions.Range("F21").Value = "'0.004616475810352206"
$wsPositions.Range("G21").Value = "'461.03"
$wsPositions.Range("F22").Value = "'0.018195551574008137"
$wsPositions.Range("G22").Value = "'4283.45"
$wsPositions.Range("F23").Value = "'0.009707932720953876"
$wsPositions.Range("G23").Value = "'11869.79"
$wsPositions.Range("F24").Value = "'0.00287701430496335"
$wsPositions.Range("G24").Value = "'3152.92"
$wsPositions.Range("F27").Value = "'0.002418018646897028"
$wsPositions.Range("G27").Value = "'4590.86"
$wsPositions.Range("F28").Value = "'0.007102039569921058"
$wsPositions.Range("G28").Value = "'2410.03"
$wsPositions.Range("F29").Value = "'0.003966212891232665"
$wsPositions.Range("G29").Value = "'1873.17"
$wsPositions.Range("F30").Value = "'0.005691268976051839"
$wsPositions.Range("G30").Value = "'2880.65"
$wsPositions.Range("F31").Value = "'0.00272531448009951"
$wsPositions.Range("G31").Value = "'1430.57"
$wsPositions.Range("F33").Value = "'0.015823245848571587"
$wsPositions.Range("G33").Value = "'12614.43"
$wsPositions.Range("F34").Value = "'0.000011114592899426"
$wsPositions.Range("G34").Value = "'12.65"
$wsPositions.Range("F35").Value = "'0.016857539829686085"
$wsPositions.Range("G35").Value = "'13988.78"
$wsPositions.Range("F36").Value = "'0.001766955158360452"
$wsPositions.Range("G36").Value = "'3093.36"
$wsPositions.Range("F37").Value = "'0.000023063182660471"
$wsPositions.Range("G37").Value = "'461.17"
$wsPositions.Range("F38").Value = "'0.002789855269261671"
$wsPositions.Range("G38").Value = "'2926.24"
$wsPositions.Range("F39").Value = "'0.001362119460166577"
$wsPositions.Range("G39").Value = "'2555.82"
$wsPositions.Range("F40").Value = "'0.001363782155033188"
$wsPositions.Range("G40").Value = "'865.82"
$wsPositions.Range("F42").Value = "'0.000024499472188667"
$wsPositions.Range("G42").Value = "'3668.12"
$wsPositions.Range("F43").Value = "'0.003774514860781313"
$wsPositions.Range("F44").Value = "'0.002694559931617619"
$wsPositions.Range("G44").Value = "'11235.91"
$wsPositions.Range("F47").Value = "'0.01213908628442389"
$wsPositions.Range("G47").Value = "'6101.05"
$wsPositions.Range("G48").Value = "'8041.67"
$wsPositions.Range("F50").Value = "'0.008911991647473406"
$wsPositions.Range("G50").Value = "'11838.38"
$wsPositions.Range("F51").Value = "'0.000974441702683132"
$wsPositions.Range("G51").Value = "'1687.82"
$wsPositions.Range("F53").Value = "'0"
$wsPositions.Range("G53").Value = "'1777.3"
$wsPositions.Range("F54").Value = "'0.000117240627169197"
$wsPositions.Range("G54").Value = "'2770.77"
$wsPositions.Range("F56").Value = "'0.002300783678480231"
$wsPositions.Range("G56").Value = "'3745.26"

# --- LP Rewards sheet: periodFeesCurrency1_formatted / reward_formatted / totalFeesCommonDenominator (columns B-E) ---
$wsLpRewards = $wb.Worksheets.Item("LP Rewards")
$wsLpRewards.Range("D3").Value = "'3405.92"
$wsLpRewards.Range("D4").Value = "'1602.81"
$wsLpRewards.Range("D5").Value = "'0.14"
$wsLpRewards.Range("D6").Value = "'620.32"
$wsLpRewards.Range("D7").Value = "'27562.06"
$wsLpRewards.Range("D8").Value = "'1316.05"
$wsLpRewards.Range("C9").Value = "'2508.68"
$wsLpRewards.Range("D9").Value = "'4438.32"
$wsLpRewards.Range("E9").Value = "'477571"
$wsLpRewards.Range("D11").Value = "'2068.83"
$wsLpRewards.Range("D12").Value = "'27734.77"
$wsLpRewards.Range("D13").Value = "'262071.32"
$wsLpRewards.Range("B14").Value = "'0.184406515513322905"
$wsLpRewards.Range("C14").Value = "'158139.53"
$wsLpRewards.Range("D14").Value = "'295216.35"
$wsLpRewards.Range("E14").Value = "'31765750"
$wsLpRewards.Range("D15").Value = "'6621.74"
$wsLpRewards.Range("D16").Value = "'8216.27"
$wsLpRewards.Range("B17").Value = "'0.005691268976051839"
$wsLpRewards.Range("C17").Value = "'2880.65"
$wsLpRewards.Range("D17").Value = "'7252.48"
$wsLpRewards.Range("E17").Value = "'780379"
$wsLpRewards.Range("B18").Value = "'0.000011114592899426"
$wsLpRewards.Range("C18").Value = "'12.65"
$wsLpRewards.Range("D18").Value = "'20.68"
$wsLpRewards.Range("E18").Value = "'2226"
